$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.881.93"
$ws.Range("E2").Value = '  -0.83%  '
$ws.Range("D3").Value = "'1.667.69"
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = "'215.56"
$ws.Range("E5").Value = '  +0.12%  '
$ws.Range("D6").Value = "'0.529"
$ws.Range("E6").Value = '  +4.09%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  +1.53%  '
$ws.Range("E9").Value = '  +0.71%  '
$ws.Range("E10").Value = '  +3.55%  '
$ws.Range("E11").Value = '  +3.97%  '
$ws.Range("D12").Value = "'1.902.15"
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").Value = "'1.629.95"
$ws.Range("E13").Value = '  -1.63%  '
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("E15").Value = '  +0.78%  '
$ws.Range("D16").Value = "'65.97"
$ws.Range("E16").Value = '  +1.59%  '
$ws.Range("D17").Value = "'26.898.91"
$ws.Range("E17").Value = '  -0.76%  '
$ws.Range("D18").Value = "'231.69"
$ws.Range("E18").Value = '  -3.79%  '
$ws.Range("D19").Value = "'7.80"
$ws.Range("E19").Value = '  -0.53%  '
$ws.Range("E20").Value = '  +0.63%  '
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("D23").Value = "'2.22"
$ws.Range("E23").Value = '  -1.69%  '
$ws.Range("D24").Value = "'9.18"
$ws.Range("E24").Value = '  -1.01%  '
$ws.Range("D25").Value = "'145.60"
$ws.Range("E25").Value = '  -0.38%  '
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("D27").Value = "'0.115"
$ws.Range("E27").Value = '  +1.00%  '
$ws.Range("E28").Value = '  +0.19%  '
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("E30").Value = '  -0.23%  '
$ws.Range("E31").Value = '  +0.34%  '
$ws.Range("E32").Value = '  +1.70%  '
$ws.Range("D33").Value = "'1.465.21"
$ws.Range("E33").Value = '  -3.82%  '
$ws.Range("E34").Value = '  +3.44%  '
$ws.Range("E35").Value = '  +3.42%  '
$ws.Range("D36").Value = "'2.41"
$ws.Range("E36").Value = '  -0.38%  '
$ws.Range("D37").Value = "'0.898"
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("D38").Value = "'0.572"
$ws.Range("E38").Value = '  -0.89%  '
$ws.Range("D39").Value = "'0.0169"
$ws.Range("E39").Value = '  -0.08%  '
$ws.Range("D40").Value = "'5.81"
$ws.Range("E40").Value = '  -2.36%  '
$ws.Range("E41").Value = '  +0.21%  '
$ws.Range("D42").Value = "'2.28"
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("E43").Value = '  +6.53%  '
$ws.Range("D44").Value = "'65.74"
$ws.Range("E44").Value = '  +0.26%  '
$ws.Range("D45").Value = "'1.812.49"
$ws.Range("E45").Value = '  +0.95%  '
$ws.Range("E46").Value = '  +0.85%  '
$ws.Range("D47").Value = "'90.28"
$ws.Range("E47").Value = '  -0.27%  '
$ws.Range("D48").Value = "'1.52"
$ws.Range("E48").Value = '  -0.90%  '
$ws.Range("E49").Value = '  +2.79%  '
$ws.Range("D50").Value = "'0.0508"
$ws.Range("E50").Value = '  +1.09%  '
$ws.Range("D51").Value = "'7.55"
$ws.Range("E51").Value = '  +0.32%  '
